$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.027.14'
$ws.Range("E2").Value = '  +8.38%  '

$ws.Range("D3").Value = '1.756.21'
$ws.Range("E3").Value = '  +4.86%  '

$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").Value = '''335.61'
$ws.Range("E5").Value = '  +1.12%  '

$ws.Range("D6").Value = '''0.9962'
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").Value = '''0.3747'
$ws.Range("E7").Value = '  +2.67%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '''48.86'
$ws.Range("E8").Value = '  +3.94%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3412'
$ws.Range("E9").Value = '  +5.46%  '

$ws.Range("D10").Value = '''1.195'
$ws.Range("E10").Value = '  +4.68%  '

$ws.Range("D11").Value = '''0.07563'
$ws.Range("E11").Value = '  +6.11%  '

$ws.Range("D12").Value = '''0.9995'
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("D13").Value = '''6.412'
$ws.Range("E13").Value = '  +5.36%  '

$ws.Range("D14").Value = '''20.64'
$ws.Range("E14").Value = '  +4.94%  '

$ws.Range("D15").Value = '''7.100'
$ws.Range("E15").Value = '  +6.93%  '

$ws.Range("D16").Value = '1.753.89'
$ws.Range("E16").Value = '  +5.09%  '

$ws.Range("E17").Value = '  +4.21%  '

$ws.Range("D18").Value = '''0.06739'
$ws.Range("E18").Value = '  +3.07%  '

$ws.Range("D19").Value = '''83.38'
$ws.Range("E19").Value = '  +5.95%  '

$ws.Range("D20").Value = '''0.9964'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").Value = '''16.93'
$ws.Range("E21").Value = '  +6.73%  '

$ws.Range("D22").Value = '''6.286'
$ws.Range("E22").Value = '  +6.24%  '

$ws.Range("D23").Value = '''12.96'
$ws.Range("E23").Value = '  +1.10%  '

$ws.Range("D24").Value = '27.024.52'
$ws.Range("E24").Value = '  +8.37%  '

$ws.Range("D25").Value = '''2.447'
$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").Value = '''1.486'
$ws.Range("E26").Value = '  +26.02%  '

$ws.Range("D27").Value = '''2.447'
$ws.Range("E27").Value = '  +2.33%  '

$ws.Range("D28").Value = '''152.17'
$ws.Range("E28").Value = '  +2.71%  '

$ws.Range("D29").Value = '''19.78'
$ws.Range("E29").Value = '  +5.72%  '

$ws.Range("D30").Value = '1.950.22'
$ws.Range("E30").Value = '  +5.20%  '

$ws.Range("D31").Value = '''133.61'
$ws.Range("E31").Value = '  +5.97%  '

$ws.Range("D32").Value = '''4.120'
$ws.Range("E32").Value = '  +0.76%  '

$ws.Range("D33").Value = '''6.078'
$ws.Range("E33").Value = '  +4.88%  '

$ws.Range("D34").Value = '''0.08652'
$ws.Range("E34").Value = '  +2.09%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '''13.02'
$ws.Range("E35").Value = '  +5.77%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '''1.691'
$ws.Range("E36").Value = '  +2.68%  '

$ws.Range("D37").Value = '''5.492'
$ws.Range("E37").Value = '  +6.42%  '

$ws.Range("D38").Value = '''0.02357'
$ws.Range("E38").Value = '  +5.45%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '''0.2203'
$ws.Range("E39").Value = '  +5.62%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '''0.06350'
$ws.Range("E40").Value = '  +5.06%  '

$ws.Range("D41").Value = '''8.604'
$ws.Range("E41").Value = '  +4.63%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.228'
$ws.Range("E42").Value = '  -0.24%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '''0.6376'
$ws.Range("E43").Value = '  +7.05%  '

$ws.Range("D44").Value = '''14.41'
$ws.Range("E44").Value = '  +4.96%  '

$ws.Range("D45").Value = '''0.9970'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").Value = '''0.6308'
$ws.Range("E46").Value = '  +10.32%  '

$ws.Range("D47").Value = '''3.861'
$ws.Range("E47").Value = '  +0.35%  '

$ws.Range("D48").Value = '''130.45'
$ws.Range("E48").Value = '  +4.84%  '

$ws.Range("D49").Value = '''2.097'
$ws.Range("E49").Value = '  +6.81%  '

$ws.Range("D50").Value = '''0.07257'
$ws.Range("E50").Value = '  +3.60%  '

$ws.Range("D51").Value = '''78.93'
$ws.Range("E51").Value = '  +5.98%  '
